# aggiornamento fino a 27/05
# Append new daily rows (14/05/2021 - 27/05/2021) to the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 256
$lastNewRow = 269

# New data: date serial (col A), nuovi pos. (col B), somma mobile 7gg. (col C),
# somma mobile 7gg. per 100mila abitanti (col D)
$data = @(
    @(44330, 1, 3, 35.34817956875221),
    @(44331, 0, 3, 35.34817956875221),
    @(44332, 0, 3, 35.34817956875221),
    @(44333, 1, 4, 47.13090609166961),
    @(44334, 0, 2, 23.5654530458348),
    @(44335, 1, 3, 35.34817956875221),
    @(44336, 0, 3, 35.34817956875221),
    @(44337, 0, 2, 23.5654530458348),
    @(44338, 2, 4, 47.13090609166961),
    @(44339, 0, 4, 47.13090609166961),
    @(44340, 0, 3, 35.34817956875221),
    @(44341, 0, 3, 35.34817956875221),
    @(44342, 0, 2, 23.5654530458348),
    @(44343, 0, 2, 23.5654530458348)
)

# Copy the formatting (date number format / style) of the last existing
# row's date cell down into the new date cells before filling in values.
$ws.Range("A255").Copy()
$ws.Range("A" + $firstNewRow + ":A" + $lastNewRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$r = $firstNewRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $r = $r + 1
}
